# Updated cryptos list on Mon May  6 14:41:49 UTC 2024 with GitHub Actions
#
# Refreshes Price (D) and Volume(1h) (E) columns with the latest scrape,
# and re-orders a couple of coins whose rank flipped (rows 17/18 and
# 40/41/42) by rewriting Coin (B) and Link (C) in place.
#
# Price cells that look like plain numbers (e.g. "590.50") need to stay
# text, matching the rest of the sheet, so those are briefly switched to
# a text NumberFormat while the value is written and then restored to
# the "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.204.59"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "3.105.96"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.42%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.102.30"
$ws.Range("E8").Value = "  -1.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.08%  "

$ws.Range("D15").Value = "3.620.72"
$ws.Range("E15").Value = "  -0.89%  "

$ws.Range("E16").Value = "  -1.91%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.978.85"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("D19").Value = "3.102.93"
$ws.Range("E19").Value = "  -1.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.737"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.86%  "

$ws.Range("E24").Value = "  +3.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.85%  "

$ws.Range("E35").Value = "  +0.67%  "

$ws.Range("E36").Value = "  -0.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.93%  "

$ws.Range("E38").Value = "  +2.13%  "

$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.16%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "455.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.42%  "

$ws.Range("E43").Value = "  +1.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0370"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").Value = "2.855.54"
$ws.Range("E45").Value = "  -1.62%  "

$ws.Range("E46").Value = "  +2.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.70%  "
